# Git007 - updated september month expenses
# Add a new "September" worksheet (after "August") with the month's
# grocery/household expense data, mirroring the layout of "August".

$wb = $excel.ActiveWorkbook
$augustSheet = $wb.Worksheets.Item("August")

# Insert the new sheet right after August.
$ws2 = $wb.Worksheets.Add($null, $augustSheet)
$ws2.Name = "September"

# Copy header row (A1:G1) and the first data row's formatting (A2:G2) from
# August so the new sheet inherits the same cell styles (bold header band,
# bordered data rows) without hand-rolling font/fill/border objects.
$augustSheet.Range("A1:G1").Copy($ws2.Range("A1:G1"))
$augustSheet.Range("A2:G2").Copy($ws2.Range("A2:G48"))

function Set-RowData($row, $brand, $product, $price, $qty, $type, $month) {
    $ws2.Cells.Item($row, 1).Value = $brand
    $ws2.Cells.Item($row, 2).Value = $product
    $ws2.Cells.Item($row, 3).Value = $price
    $ws2.Cells.Item($row, 4).Value = $qty
    $ws2.Cells.Item($row, 5).Formula = "=C" + $row + "*D" + $row
    $ws2.Cells.Item($row, 6).Value = $type
    $ws2.Cells.Item($row, 7).Value = $month
}

Set-RowData 2 "Idhayam" "Nalennai 1L" 239 4 "OF" "September"
Set-RowData 3 "Maggie" "Noodles 560g" 71 1 "OF" "September"
Set-RowData 4 "Anil" "Semiya 900g" 54.5 1 "OF" "September"
Set-RowData 5 "Chakra Gold" "Tea 500g" 200.5 1 "OF" "September"
Set-RowData 6 "Surf excel" "Washing liquid 500ml" 226 2 "SG" "September"
Set-RowData 7 "Harpic cit" "Toilet Flush 500ml" 66 2 "SG" "September"
Set-RowData 8 "Harpic pow" "Toilet Flush 500ml" 65 2 "SG" "September"
Set-RowData 9 "Lizol lav" "Toilet Floor 500ml" 64 1 "SG" "September"
Set-RowData 10 "Lizol cit" "Toilet Floor 500ml" 65 1 "SG" "September"
Set-RowData 11 "Vim" "Dish wash Liquid 500ml" 88 1 "SG" "September"
Set-RowData 12 "Aachi" "Pepper powder 50g" 79 1 "NS" "September"
Set-RowData 13 "BP" "Corn flour 100g" 22.5 1 "GF" "September"
Set-RowData 14 "LG" "Perungayam Katti 50g" 41 1 "NS" "September"
Set-RowData 15 "Sakthi" "Malli Thool 500g" 128 2 "NS" "September"
Set-RowData 16 "Ponnu" "Seeragam 50g" 18 1 "NS" "September"
Set-RowData 17 "LG" "Perungayam 100g" 85 1 "NS" "September"
Set-RowData 18 "Ponnu" "Black Mookadalai 1kg" 77 1 "GF" "September"
Set-RowData 19 "Sakthi" "ManjalThool 50g" 8 1 "NS" "September"
Set-RowData 20 "Udhayam" "Ghee 500ml" 240.5 1 "OF" "September"
Set-RowData 21 "Gokul Santol" "Powder 70gm" 48.5 1 "SG" "September"
Set-RowData 22 "Twinkle" "Scrubber pad" 9.5 3 "SG" "September"
Set-RowData 23 "SB" "Scrubber" 14 2 "SG" "September"
Set-RowData 24 "Ponnu" "Kanaramani 500g" 72.5 1 "GF" "September"
Set-RowData 25 "Dove" "Bath Soap 75g" 43.5 1 "SG" "September"
Set-RowData 26 "Ponnu" "Pattai 50g" 20 1 "NS" "September"
Set-RowData 27 "Ponnu" "Cloves 50g" 105 1 "NS" "September"
Set-RowData 28 "Ponnu" "Kitchen Towel" 132 1 "SG" "September"
Set-RowData 29 "Cif" "Dishwash liquid 250ml" 59.5 1 "SG" "September"
Set-RowData 30 "Lifeboy" "Hand Wash 900ml(sache)" 150 1 "SG" "September"
Set-RowData 31 "Yardly" "Body Deodarant" 178 1 "SG" "September"
Set-RowData 32 "Hamam" "Bath Soap 100g" 25 4 "SG" "September"
Set-RowData 33 "Ponnu" "Samba Ravai 500g" 50 1 "OF" "September"
Set-RowData 34 "Ponnu" "Cashewnut 200g" 180 1 "NS" "September"
Set-RowData 35 "Ponnu" "Cumin seeds 200g" 50 1 "NS" "September"
Set-RowData 36 "Ponnu" "Pepper 200g" 216 1 "NS" "September"
Set-RowData 37 "Ponnu" "Mustard 200g" 20 2 "NS" "September"
Set-RowData 38 "Ponnu" "Green Moong dhal" 57 1 "GF" "September"
Set-RowData 39 "Ponnu" "Aval 500g" 32.5 1 "OF" "September"
Set-RowData 40 "Ponnu" "Raisins 500g" 195 1 "NS" "September"
Set-RowData 41 "Ponnu" "Green peas 200g" 16 1 "GF" "September"
Set-RowData 42 "Ponnu" "Poppy 50g" 32 1 "NS" "September"
Set-RowData 43 "Ponnu" "Methi 200g" 19 1 "NS" "September"
Set-RowData 44 "Ponnu" "Cardamom 25g" 42 1 "NS" "September"
Set-RowData 45 "Anil" "Roasted Rava 1kg" 50 1 "OF" "September"
Set-RowData 46 "Pepsodent" "Toothpaste 80g" 42.5 1 "SG" "September"
Set-RowData 47 "Ponnu" "Kalpasi 10g" 10 1 "NS" "September"
Set-RowData 48 "Ponnu" "Aacturm powder 100g" 15 2 "NS" "September"

# Column widths (best-fit-like) matching the committed sheet.
$ws2.Columns.Item(1).ColumnWidth = 11.5
$ws2.Columns.Item(2).ColumnWidth = 22.67
$ws2.Columns.Item(3).ColumnWidth = 11.67
$ws2.Columns.Item(4).ColumnWidth = 10.5
$ws2.Columns.Item(5).ColumnWidth = 9.67
$ws2.Columns.Item(6).ColumnWidth = 5.83
$ws2.Columns.Item(7).ColumnWidth = 10

# Selection / active-sheet bookkeeping: September becomes the active tab
# with G3:G48 selected, while August loses its old selection/scroll state
# and settles on A2.
$augustSheet.Range("A2").Select()
$ws2.Activate()
$ws2.Range("G3:G48").Select()
